# Commit message: "Add files via upload"
# Content-level changes in the diff:
#   1. Cell B2 (sheet "ImgList") text value changed: "alilo888" -> "asma"
#   2. Active selection on the sheet changed from A3 to E7
# (The remaining diff hunks are Excel-generated boilerplate/noise: new
#  xr/xr2/xr3 namespace declarations, a regenerated revisionPtr GUID, a
#  regenerated xr2:uid, and the author's local absPath — none of these are
#  reachable/meaningful via the Excel object model, so they're left alone.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "asma"

[void]$ws.Range("E7").Select()
